# Update the fund name used in the capital calls sample upload.
# "Agri Fund" -> "SAAS Fund" for all data rows in column A,
# and move the active selection to A4 (as captured by the author's edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# Replace every occurrence of "Agri Fund" in column A with "SAAS Fund".
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "Agri Fund") {
        $cell.Value = "SAAS Fund"
    }
}

# Update the selection to match the saved view state (A4).
$ws.Range("A4").Select()
